# Update "想去人数" (F) and "最低票价" (G) figures across sheets,
# reflecting the refreshed scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 129
$ws1.Range("F5").Value = 5485
$ws1.Range("F7").Value = 4
$ws1.Range("F10").Value = 2438
$ws1.Range("F11").Value = 80
$ws1.Range("F12").Value = 69
$ws1.Range("G12").Value = 70
$ws1.Range("F14").Value = 2288
$ws1.Range("F15").Value = 180

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 98

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 129
$ws4.Range("F5").Value = 5485
$ws4.Range("F6").Value = 98
$ws4.Range("F8").Value = 4
$ws4.Range("F12").Value = 2438
$ws4.Range("F13").Value = 80
$ws4.Range("F14").Value = 69
$ws4.Range("G14").Value = 70
$ws4.Range("F17").Value = 2288
$ws4.Range("F18").Value = 180

$wb.Save()
